# Fix regtime/borntime data so that regtime (column D) is always after
# borntime (column E) for every reader record.
#
# For each row:
#   1. Read the existing D (regtime) and E (borntime) serial date values.
#   2. Round each value to the nearest whole second.
#   3. Try to "fix" each rounded value by swapping its day and month
#      components (this mirrors the original data-generation bug where
#      day/month got transposed for some records). If the swap produces an
#      invalid date (day > 12) the value is instead written out as literal
#      text in dd/MM/yyyy HH:mm:ss format (mirroring what the row's
#      original generator script produced when it could not re-interpret
#      the date).
#   4. If the original regtime (D) was already after the original borntime
#      (E), each column keeps its own corrected value. Otherwise the two
#      corrected values are swapped between the columns so that regtime
#      ends up after borntime.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function MakeDate($year, $month, $day, $hour, $minute, $second) {
    $base = [datetime]::FromOADate(2)   # 01/01/1900
    return $base.AddYears($year - 1900).AddMonths($month - 1).AddDays($day - 1).AddHours($hour).AddMinutes($minute).AddSeconds($second)
}

function RoundToSecondOADate($v) {
    $totalSeconds = $v * 86400.0
    $roundedSeconds = [Math]::Floor($totalSeconds + 0.5)
    return $roundedSeconds / 86400.0
}

# Returns a 2-element array: @("num", <double>) or @("str", <text>)
function TransformSerial($serial) {
    $rounded = RoundToSecondOADate $serial
    $dt = [datetime]::FromOADate($rounded)
    $year = $dt.Year
    $month = $dt.Month
    $day = $dt.Day
    $hour = $dt.Hour
    $minute = $dt.Minute
    $second = $dt.Second

    if ($day -ge 1 -and $day -le 12) {
        $swapped = MakeDate $year $day $month $hour $minute $second
        return @("num", $swapped.ToOADate())
    } else {
        return @("str", $dt.ToString("dd/MM/yyyy HH:mm:ss"))
    }
}

function ApplyResult($cell, $result) {
    if ($result[0] -eq "num") {
        $cell.Value = $result[1]
    } else {
        $cell.Value = $result[1]
    }
}

$lastRow = 119
for ($row = 2; $row -le $lastRow; $row++) {
    $dCell = $ws.Cells.Item($row, 4)
    $eCell = $ws.Cells.Item($row, 5)

    $dOld = $dCell.Value2
    $eOld = $eCell.Value2

    $dRes = TransformSerial $dOld
    $eRes = TransformSerial $eOld

    if ($dOld -gt $eOld) {
        # regtime already after borntime: transform each value in place
        ApplyResult $dCell $dRes
        ApplyResult $eCell $eRes
    } else {
        # regtime was before borntime: swap the two corrected values
        ApplyResult $dCell $eRes
        ApplyResult $eCell $dRes
    }
}
